# Apply the updates to the "Liquidación" worksheet:
#  - C17: change from formula (=Datos!J17) to a literal value (10)
#  - C18: literal value 6 -> 8
#  - J18: literal value 0 -> 3
#  - C21: date serial 44932 -> 45194  (2023-01-06 -> 2023-09-25)
#  - C22: date serial 44927 -> 45108  (2023-01-01 -> 2023-07-01)
#
# All other changed cells in the workbook are formula results that
# automatically recompute once these inputs change (DAYS360 calculations,
# prorated liquidation totals, the "amount in words" helper sheet, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Liquidación")

# Disconnect C17 from the Datos!J17 formula and set the literal value.
$ws.Range("C17").Value = 10

# Update the days-worked / pending counters.
$ws.Range("C18").Value = 8
$ws.Range("J18").Value = 3

# Update the vacation / interest cut-off dates.
$ws.Range("C21").Value = 45194
$ws.Range("C22").Value = 45108

# Force a full recalculation so every dependent formula (DAYS360 totals,
# liquidation summary, and the amount-in-words conversion sheet) reflects
# the new inputs.
$excel.CalculateFullRebuild()

# Restore the sheet/selection state recorded after the edits were made.
$ws.Activate()
$ws.Range("J36").Select() | Out-Null
